$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete row 44 entirely; remaining rows (45-73) shift up to become 44-72
$ws.Rows.Item(44).Delete()

# Update selection to reflect where the cursor ended up after the delete (A44)
$ws.Activate()
$ws.Range("A44").Select()
